# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the 423dd735 file's
# xlf round-trip (zh-cn + de-de), and roll the Overview sheet's
# "Latest HO Xliff Generate Date" for 423dd735 forward to match.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-13 00:59:42"
$wsZhCn.Range("K2").Value = "2016-08-13 01:00:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-13 00:59:49"
$wsDeDe.Range("K2").Value = "2016-08-13 01:00:29"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-13 00:59:49"
